{"js": "// Remove the \"[WIP] \" prefix from the four \"Current templates\" list items\n// in the README: \"Report\", \"eBook\", \"Thesis\", \"Website\" are no longer WIP.\nconst body = context.document.body;\n\nconst replacements = [\n  [\"[WIP] Report\", \"Report\"],\n  [\"[WIP] eBook\", \"eBook\"],\n  [\"[WIP] Thesis\", \"Thesis\"],\n  [\"[WIP] Website\", \"Website\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const result of results.items) {\n    result.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the \"[WIP] \" prefix from the four \"Current templates\" list items\n# in the README: \"Report\", \"eBook\", \"Thesis\", \"Website\" are no longer WIP.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"[WIP] Report\", \"Report\"),\n    @(\"[WIP] eBook\", \"eBook\"),\n    @(\"[WIP] Thesis\", \"Thesis\"),\n    @(\"[WIP] Website\", \"Website\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $false, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)  # wdReplaceAll\n}\n"}
